$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "AddCustomerTest"

# Populate cells in the same order the original authoring tool wrote them
# (this determines the shared-strings table order).
$ws.Range("A1").Value = "firstName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "postCode"
$ws.Range("A2").Value = "Anne"
$ws.Range("B2").Value = "Zimmermann"
$ws.Range("C2").Value = "89046-305"
$ws.Range("D1").Value = "alertText"
$ws.Range("D2").Value = "Customer added successfully"

# Header row styling: JetBrains Mono 10pt, vertically centered
$header = $ws.Range("A1:D1")
$header.Font.Name = "JetBrains Mono"
$header.Font.Size = 10
$header.Font.Color = 13023145
$header.VerticalAlignment = -4108

# Column widths (characters)
$ws.Columns.Item(1).ColumnWidth = 9.799479166666666
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 8.709635416666666
$ws.Columns.Item(4).ColumnWidth = 38.983072916666664

# Selection
$ws.Range("E3").Select()
